$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - copy formatting from existing header (G1)
# so it matches the other bold/bordered/centered header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New Save value for the data row
$ws.Range("H2").Value = 1
